$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 (entered first, completely) ---
$ws.Range("A14").Value = "New Person"
$ws.Range("B14").Value = "New position"
$ws.Range("C14").Value = 780223068
$ws.Range("D14").Value = "newperson@mail.com"

# --- Column A for the remaining new rows ---
$ws.Range("A15").Value = "New BPerson"
$ws.Range("A16").Value = "New C Person"
$ws.Range("A17").Value = "New D Person"
$ws.Range("A18").Value = "New E Person"
$ws.Range("A19").Value = "New F Person"

# --- Column B (position, same text reused) and C (phone numbers) ---
$ws.Range("B15").Value = "New position"
$ws.Range("B16").Value = "New position"
$ws.Range("B17").Value = "New position"
$ws.Range("B18").Value = "New position"
$ws.Range("B19").Value = "New position"

$ws.Range("C15").Value = 780223069
$ws.Range("C16").Value = 780223070
$ws.Range("C17").Value = 780223071
$ws.Range("C18").Value = 780223072
$ws.Range("C19").Value = 780223073

# --- Column D (emails) for the remaining new rows ---
$ws.Range("D15").Value = "newAperson@mail.com"
$ws.Range("D16").Value = "newOPerson@mail.com"
$ws.Range("D17").Value = "person@mail.com"
$ws.Range("D18").Value = "newemali@mail.com"
$ws.Range("D19").Value = "testing@mail.com"

# --- Hyperlinks (mailto:) for each new email cell ---
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:newperson@mail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:newAperson@mail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D16"), "mailto:newOPerson@mail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D17"), "mailto:person@mail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D18"), "mailto:newemali@mail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D19"), "mailto:testing@mail.com") | Out-Null

# Re-apply the existing "Hyperlink" cell style to the new email cells so they
# carry the same style index as the rest of the email column.
$ws.Range("D14:D19").Style = "Hyperlink"

# --- Selection, as captured in the saved workbook ---
$ws.Range("G17").Select() | Out-Null
